$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "62.973.02"
$ws.Range("E2").Value = "  +4.95%  "
$ws.Range("D3").Value = "2.454.52"
$ws.Range("E3").Value = "  +3.24%  "
$ws.Range("E4").Value = "  +0.12%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "574.05"
$ws.Range("D5").NumberFormat = "General"
$ws.Range("E5").Value = "  +2.45%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "146.13"
$ws.Range("D6").NumberFormat = "General"
$ws.Range("E6").Value = "  +5.74%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.00"
$ws.Range("D7").NumberFormat = "General"
$ws.Range("E7").Value = "  -0.08%  "
$ws.Range("E8").Value = "  +2.11%  "
$ws.Range("D9").Value = "2.453.36"
$ws.Range("E9").Value = "  +3.27%  "
$ws.Range("E10").Value = "  +5.23%  "
$ws.Range("E12").Value = "  +2.79%  "
$ws.Range("E13").Value = "  +4.22%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "27.31"
$ws.Range("D14").NumberFormat = "General"
$ws.Range("E14").Value = "  +6.24%  "
$ws.Range("E15").Value = "  +7.32%  "
$ws.Range("D16").Value = "2.902.70"
$ws.Range("E16").Value = "  +3.88%  "
$ws.Range("D17").Value = "62.850.44"
$ws.Range("E17").Value = "  +5.00%  "
$ws.Range("D18").Value = "2.440.40"
$ws.Range("E18").Value = "  +2.71%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.93"
$ws.Range("D19").NumberFormat = "General"
$ws.Range("E19").Value = "  -2.01%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "10.97"
$ws.Range("D20").NumberFormat = "General"
$ws.Range("E20").Value = "  +4.24%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "328.70"
$ws.Range("D21").NumberFormat = "General"
$ws.Range("E21").Value = "  +2.32%  "
$ws.Range("E22").Value = "  +1.97%  "
$ws.Range("E23").Value = "  +12.90%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.00"
$ws.Range("D24").NumberFormat = "General"
$ws.Range("E24").Value = "  -0.19%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "65.69"
$ws.Range("D25").NumberFormat = "General"
$ws.Range("E25").Value = "  +2.59%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "635.77"
$ws.Range("D26").NumberFormat = "General"
$ws.Range("E26").Value = "  +13.92%  "
$ws.Range("B27").Value = "Binance-PegBSC-USD"
$ws.Range("C27").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.10"
$ws.Range("D27").NumberFormat = "General"
$ws.Range("E27").Value = "  +10.26%  "
$ws.Range("B28").Value = "Aptos"
$ws.Range("C28").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "8.51"
$ws.Range("D28").NumberFormat = "General"
$ws.Range("E28").Value = "  +4.27%  "
$ws.Range("B29").Value = "PEPE"
$ws.Range("C29").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D29").Value = "0.0₃0989"
$ws.Range("E29").Value = "  +6.71%  "
$ws.Range("B30").Value = "WrappedeETH"
$ws.Range("C30").Value = "https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth"
$ws.Range("D30").Value = "2.534.16"
$ws.Range("E30").Value = "  +1.85%  "
$ws.Range("B31").Value = "InternetComputer(DFINITY)"
$ws.Range("C31").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "8.21"
$ws.Range("D31").NumberFormat = "General"
$ws.Range("E31").Value = "  +1.93%  "
$ws.Range("B32").Value = "Fetch.AI"
$ws.Range("C32").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.42"
$ws.Range("D32").NumberFormat = "General"
$ws.Range("E32").Value = "  +8.91%  "
$ws.Range("B33").Value = "Kaspa"
$ws.Range("C33").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.139"
$ws.Range("D33").NumberFormat = "General"
$ws.Range("E33").Value = "  +5.94%  "
$ws.Range("B34").Value = "PancakeSwap"
$ws.Range("C34").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.86"
$ws.Range("D34").NumberFormat = "General"
$ws.Range("E34").Value = "  +4.17%  "
$ws.Range("B35").Value = "ImmutableX"
$ws.Range("C35").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.50"
$ws.Range("D35").NumberFormat = "General"
$ws.Range("E35").Value = "  +4.46%  "
$ws.Range("B36").Value = "FirstDigitalUSD"
$ws.Range("C36").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.999"
$ws.Range("D36").NumberFormat = "General"
$ws.Range("E36").Value = "  +0.02%  "
$ws.Range("B37").Value = "NEARProtocol"
$ws.Range("C37").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "4.76"
$ws.Range("D37").NumberFormat = "General"
$ws.Range("E37").Value = "  +4.59%  "
$ws.Range("B38").Value = "PolygonEcosystemToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.375"
$ws.Range("D38").NumberFormat = "General"
$ws.Range("E38").Value = "  +2.22%  "
$ws.Range("B39").Value = "Monero"
$ws.Range("C39").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "153.03"
$ws.Range("D39").NumberFormat = "General"
$ws.Range("E39").Value = "  -0.24%  "
$ws.Range("B40").Value = "RenderToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/vfo5XYwcV+rendertoken-render"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "5.41"
$ws.Range("D40").NumberFormat = "General"
$ws.Range("E40").Value = "  +8.38%  "
$ws.Range("B41").Value = "EthereumClassic"
$ws.Range("C41").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "18.69"
$ws.Range("D41").NumberFormat = "General"
$ws.Range("E41").Value = "  +3.02%  "
$ws.Range("B42").Value = "dogwifhat"
$ws.Range("C42").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.79"
$ws.Range("D42").NumberFormat = "General"
$ws.Range("E42").Value = "  +15.14%  "
$ws.Range("B43").Value = "Stacks"
$ws.Range("C43").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.77"
$ws.Range("D43").NumberFormat = "General"
$ws.Range("E43").Value = "  +7.57%  "
$ws.Range("B44").Value = "USDe"
$ws.Range("C44").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.999"
$ws.Range("D44").NumberFormat = "General"
$ws.Range("E44").Value = "  +0.01%  "
$ws.Range("B45").Value = "BabyDogeCoin"
$ws.Range("C45").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D45").Value = "0.0₆0284"
$ws.Range("E45").Value = "  -5.09%  "
$ws.Range("B46").Value = "Aave"
$ws.Range("C46").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "145.30"
$ws.Range("D46").NumberFormat = "General"
$ws.Range("E46").Value = "  +3.75%  "
$ws.Range("B47").Value = "Filecoin"
$ws.Range("C47").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.60"
$ws.Range("D47").NumberFormat = "General"
$ws.Range("E47").Value = "  +2.37%  "
$ws.Range("B48").Value = "InjectiveProtocol"
$ws.Range("C48").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "20.52"
$ws.Range("D48").NumberFormat = "General"
$ws.Range("E48").Value = "  +7.49%  "
$ws.Range("B49").Value = "Mantle"
$ws.Range("C49").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.602"
$ws.Range("D49").NumberFormat = "General"
$ws.Range("E49").Value = "  +2.86%  "
$ws.Range("B50").Value = "Hedera"
$ws.Range("C50").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0517"
$ws.Range("D50").NumberFormat = "General"
$ws.Range("E50").Value = "  +3.33%  "
$ws.Range("B51").Value = "WhiteBITCoin"
$ws.Range("C51").Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "13.04"
$ws.Range("D51").NumberFormat = "General"
$ws.Range("E51").Value = "  +11.54%  "
